$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.516.37"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.239.62"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.857"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "2.256.16"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "42.299.97"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000105"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.06%  "
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +45.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.02%  "
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0816"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.38%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  +3.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.62%  "
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -4.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.88%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.96%  "
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +5.42%  "
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
